$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OCT-2020")
$ws.Activate()

# Row 28: B28 was the mistyped text "10/126/2020" -> replace with the real
# date value for 10/26/2020 (Excel serial 44130). This also removes the
# now-unused shared string, shrinking the shared-strings table.
$ws.Range("B28").Value = 44130

# Rows 29-31: shift each entry's date forward by one day (these rows had
# all mistakenly been stamped with the same date).
$ws.Range("B29").Value = 44131
$ws.Range("B30").Value = 44132
$ws.Range("B31").Value = 44133

# Add a new row 32 (entry #30) for 10/30/2020, copying the formatting from
# row 31 (the row immediately above) since it is identical other than the
# serial number and date.
$ws.Range("A31:G31").Copy($ws.Range("A32:G32"))
$ws.Range("A32").Value = 30
$ws.Range("B32").Value = 44134

$ws.Range("C32").Value = "QMVAR"
$ws.Range("D32").Value = "issue fixing"
$ws.Range("F32").Value = "WIP"

# Update the view so the window shows the newly-added rows, matching the
# saved selection/scroll state.
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C39").Select()
